$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (D1 becomes Substitutions, subsequent headers shift right)
$ws.Range("D1").Value = "Substitutions"
$ws.Range("E1").Value = "StudyType"
$ws.Range("F1").Value = "Clade"
$ws.Range("G1").Value = "ResistanceLevel"
$ws.Range("H1").Value = "FoldChange"

# Fill the new clinical-record row (row 2) with new values first, reused
# columns (C2/G2) left to adopt existing shared strings, then A2 last.
$ws.Range("B2").Value = "Lenz_et_al_2013"
$ws.Range("D2").Value = "168Q"
$ws.Range("E2").Value = "clinical"
$ws.Range("F2").Value = "3a"
$ws.Range("C2").Value = "NS3"
$ws.Range("G2").Value = "resistant"
$ws.Range("H2").Value = ""
$ws.Range("A2").Value = "simeprevir"

# Add the in_vitro record as a new row 3
$ws.Range("A3").Value = "simeprevir"
$ws.Range("B3").Value = "Lenz_et_al_2013"
$ws.Range("C3").Value = "NS3"
$ws.Range("D3").Value = "168Q"
$ws.Range("E3").Value = "in_vitro"
$ws.Range("F3").Value = "1b"
$ws.Range("G3").Value = "resistant"
$ws.Range("H3").Value = 385

# New SVR column added last
$ws.Range("I1").Value = "SVR"

# A hyperlink was experimented with (and removed) during editing; Excel
# keeps the Hyperlink/Followed Hyperlink style definitions it minted even
# after the link itself is gone.
$ws.Hyperlinks.Add($ws.Range("Z1"), "", "A1", [Type]::Missing, "temp") | Out-Null
$ws.Hyperlinks.Delete()
$ws.Range("Z1").Clear()

# Update view/selection state
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("B6").Select()
